$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Personne" column (C) values — répartition des tâches
$ws.Range("C1").Value = "Personne"
$ws.Range("C2").Value = "Evan"
$ws.Range("C3").Value = "Evan"
$ws.Range("C4").Value = "Armand"
$ws.Range("C5").Value = "Armand"
$ws.Range("C6").Value = "Evan"
$ws.Range("C7").Value = "Armand / Evan"
$ws.Range("C8").Value = "Evan"
$ws.Range("C9").Value = "Evan / Armand"
$ws.Range("C10").Value = "Evan / Armand"

# Match formatting: C1 header like A1 (bold, centered), C2:C10 like A2:A10 (centered)
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

$ws.Range("A2:A10").Copy()
$ws.Range("C2:C10").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Update the view: move the active selection to E4
$ws.Range("E4").Select()
